$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.1219627926616031, 0.65096565921388, 2.0184203984331255, 3.747096004287215)"
$ws.Range("C2").Value = "NIG(0.9443484853018659, 0.6941414460604605, 5.485183473760234, 5.174350726331159)"
$ws.Range("D2").Value = "NIG(1.6449133221328798, 1.0740720095617347, 1.198802096504752, 4.1786729319112945)"
$ws.Range("E2").Value = "JSU(-1.2242333315861627, 1.1428836722772902, 2.53599059603127, 4.596828168728434)"
